$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '55.066.17'
Set-TextValue 'E2' '  -2.18%  '

Set-TextValue 'D3' '2.350.98'
Set-TextValue 'E3' '  -5.13%  '

Set-TextValue 'E4' '  -0.07%  '

Set-TextValue 'D5' '476.36'
Set-TextValue 'E5' '  -2.08%  '

Set-TextValue 'D6' '145.21'
Set-TextValue 'E6' '  -1.07%  '

Set-TextValue 'D7' '0.998'
Set-TextValue 'E7' '  +0.06%  '

Set-TextValue 'D8' '0.597'
Set-TextValue 'E8' '  +17.15%  '

Set-TextValue 'D9' '2.356.73'
Set-TextValue 'E9' '  -5.34%  '

Set-TextValue 'D10' '0.0970'
Set-TextValue 'E10' '  +0.59%  '

Set-TextValue 'D11' '5.43'
Set-TextValue 'E11' '  -5.75%  '

Set-TextValue 'D12' '0.324'
Set-TextValue 'E12' '  -1.78%  '

Set-TextValue 'D13' '0.125'
Set-TextValue 'E13' '  +1.17%  '

Set-TextValue 'D14' '2.760.23'
Set-TextValue 'E14' '  -4.94%  '

Set-TextValue 'D15' '54.984.40'
Set-TextValue 'E15' '  -2.36%  '

Set-TextValue 'D16' '20.10'
Set-TextValue 'E16' '  -4.62%  '

Set-TextValue 'D17' '0.0000131'
Set-TextValue 'E17' '  -3.26%  '

Set-TextValue 'D18' '2.349.60'
Set-TextValue 'E18' '  -5.19%  '

Set-TextValue 'D19' '4.59'
Set-TextValue 'E19' '  +1.99%  '

Set-TextValue 'D20' '317.15'
Set-TextValue 'E20' '  +0.00%  '

Set-TextValue 'D21' '9.66'
Set-TextValue 'E21' '  -3.54%  '

Set-TextValue 'D22' '1.00'
Set-TextValue 'E22' '  -0.05%  '

Set-TextValue 'D23' '5.64'
Set-TextValue 'E23' '  -2.38%  '

Set-TextValue 'D24' '56.79'
Set-TextValue 'E24' '  -2.66%  '

Set-TextValue 'D25' '0.999'
Set-TextValue 'E25' '  +0.08%  '

Set-TextValue 'D26' '0.395'
Set-TextValue 'E26' '  -3.67%  '

Set-TextValue 'D27' '0.154'
Set-TextValue 'E27' '  -4.83%  '

Set-TextValue 'D28' '2.445.09'
Set-TextValue 'E28' '  -5.52%  '

Set-TextValue 'D29' '7.18'
Set-TextValue 'E29' '  -5.22%  '

Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.11%  '

Set-TextValue 'D31' '0.0₃0759'
Set-TextValue 'E31' '  -3.98%  '

Set-TextValue 'D32' '147.22'
Set-TextValue 'E32' '  -1.36%  '

Set-TextValue 'D33' '18.14'
Set-TextValue 'E33' '  +0.17%  '

Set-TextValue 'D34' '1.47'
Set-TextValue 'E34' '  -1.74%  '

Set-TextValue 'D35' '5.07'
Set-TextValue 'E35' '  -2.27%  '

Set-TextValue 'B36' 'NEARProtocol'
Set-TextValue 'C36' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D36' '3.59'
Set-TextValue 'E36' '  -3.36%  '

Set-TextValue 'B37' 'ImmutableX'
Set-TextValue 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.09'
Set-TextValue 'E37' '  -4.27%  '

Set-TextValue 'D38' '0.820'
Set-TextValue 'E38' '  -4.60%  '

Set-TextValue 'B39' 'OKB'
Set-TextValue 'C39' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D39' '33.62'
Set-TextValue 'E39' '  -1.49%  '

Set-TextValue 'B40' 'FirstDigitalUSD'
Set-TextValue 'C40' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D40' '0.998'
Set-TextValue 'E40' '  +0.32%  '

Set-TextValue 'B41' 'Stellar'
Set-TextValue 'C41' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D41' '0.0990'
Set-TextValue 'E41' '  +7.08%  '

Set-TextValue 'B42' 'Stacks'
Set-TextValue 'C42' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D42' '1.35'
Set-TextValue 'E42' '  +1.93%  '

Set-TextValue 'B43' 'Filecoin'
Set-TextValue 'C43' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D43' '3.39'
Set-TextValue 'E43' '  -3.30%  '

Set-TextValue 'D44' '0.577'
Set-TextValue 'E44' '  -5.26%  '

Set-TextValue 'D45' '0.0525'
Set-TextValue 'E45' '  -5.66%  '

Set-TextValue 'D46' '10.14'
Set-TextValue 'E46' '  -0.35%  '

Set-TextValue 'D47' '251.95'
Set-TextValue 'E47' '  -2.76%  '

Set-TextValue 'D48' '0.0222'
Set-TextValue 'E48' '  -2.77%  '

Set-TextValue 'D49' '4.43'
Set-TextValue 'E49' '  -7.41%  '

Set-TextValue 'D50' '16.83'
Set-TextValue 'E50' '  -4.06%  '

Set-TextValue 'D51' '1.783.20'
Set-TextValue 'E51' '  -5.24%  '

